$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1): prefix with "label_"
$ws.Range("A1").Value = "label_day_1"
$ws.Range("B1").Value = "label_day_2"
$ws.Range("C1").Value = "label_day_3"
$ws.Range("D1").Value = "label_day_4"
$ws.Range("E1").Value = "label_day_5"

# Update data values (row 2)
$ws.Range("A2").Value = 0.9185514946539816
$ws.Range("B2").Value = 0.9185514946539816
$ws.Range("C2").Value = 0.9266349030356956
$ws.Range("D2").Value = 0.9346621100357814
$ws.Range("E2").Value = 0.9426341006466684
